$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Filip Bernevec, Abbas " -> "Filip " / "Bernevec" / ", Abbas " (3 runs)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Filip Bernevec, Abbas ")
$start = $r.Start

$splitA = $start + 6    # after "Filip "
$splitB = $start + 14   # after "Bernevec"

# Split off "Filip " into its own run (re-assigning FormattedText forces the
# run boundary while preserving every run-property, incl. w:cs).
$rA = $d.Range($start, $splitA)
$rA.FormattedText = $rA.FormattedText

# Split off "Bernevec" into its own run.
$rB = $d.Range($splitA, $splitB)
$rB.FormattedText = $rB.FormattedText

# ---------------------------------------------------------------------------
# 2) "Evaluate the progress of every team member in the allocated tasks"
#    -> "Evaluate the progress of " / "every team member in the allocated
#    tasks", with the (relocated) _GoBack bookmark between the two runs.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Evaluate the progress of every team member in the allocated tasks")
$start2 = $r2.Start
$splitC = $start2 + 25  # after "Evaluate the progress of "

$rC = $d.Range($start2, $splitC)
$rC.FormattedText = $rC.FormattedText

# ---------------------------------------------------------------------------
# 3) Move the _GoBack bookmark from the end of the document to the split
#    point created above.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bmRange = $d.Range($splitC, $splitC)
$d.Bookmarks.Add("_GoBack", $bmRange)
